# Reverse the "Periodo Mora" values in column E (rows 16-25) so the most
# recent period (2110) is listed first and the oldest (2101) last.
# Elimina EC anteriores y se agregan nuevos, se modifica base de datos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @(2110, 2109, 2108, 2107, 2106, 2105, 2104, 2103, 2102, 2101)

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = [string]$periodos[$i]
}
